$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# Sheet1 changes
$ws1.Range("C2").Value = "approach"
$ws1.Range("C1").Borders.Item(3).LineStyle = 1  # top
$ws1.Range("C1").Borders.Item(4).LineStyle = 1  # bottom
$ws1.Range("D1").Borders.Item(3).LineStyle = 1
$ws1.Range("D1").Borders.Item(2).LineStyle = 1  # right
$ws1.Range("D1").Borders.Item(4).LineStyle = 1

# Sheet2 changes
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"
$ws2.Range("C1").Borders.Item(3).LineStyle = 1
$ws2.Range("C1").Borders.Item(4).LineStyle = 1
$ws2.Range("D1").Borders.Item(3).LineStyle = 1
$ws2.Range("D1").Borders.Item(2).LineStyle = 1
$ws2.Range("D1").Borders.Item(4).LineStyle = 1
$ws2.Range("F1").Borders.Item(3).LineStyle = 1
$ws2.Range("F1").Borders.Item(4).LineStyle = 1
$ws2.Range("G1").Borders.Item(3).LineStyle = 1
$ws2.Range("G1").Borders.Item(2).LineStyle = 1
$ws2.Range("G1").Borders.Item(4).LineStyle = 1

# Remove G5 empty inline string cell
$ws2.Range("G5").ClearContents()
